# EZ-4133 imail templates update again
# "EZBob Ltd" -> "EZBob Ltd. (formerly known as Orange Money Ltd.)"
# in the "Notice to Guarantor" template's signature-block paragraph
# (the one that also carries the _GoBack bookmark).

$d = $word.ActiveDocument

# Locate the paragraph whose whole text is exactly "EZBob Ltd" (the
# short signature line, not the earlier "...between EZBob Ltd. (trading
# as ezbob)..." sentence).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq "EZBob Ltd") {
        $target = $para
        break
    }
}

$paraStart = $target.Range.Start

# Insert the new text right after "EZBob" (5 characters in), i.e. at the
# exact spot where the _GoBack bookmark currently sits. This pushes the
# bookmark (and the remaining old " Ltd" run) to the right, matching how
# Word relocates a collapsed bookmark when text is typed at its location.
$insertionPoint = $d.Range($paraStart + 5, $paraStart + 5)
$insertionPoint.InsertAfter(" Ltd. (form")

# The _GoBack bookmark now sits right after "(form". Grab its position so
# we can locate the leftover original " Ltd" run that still needs to
# become "erly known as Orange Money Ltd.)".
$bookmark = $d.Bookmarks.Item("_GoBack")

$paraRange = $target.Range
$tail = $d.Range($bookmark.End, $paraRange.End - 1)
$tail.Text = "erly known as Orange Money Ltd.)"
